# Apply the updated cryptocurrency Price / Volume(1h) figures (and the two
# coin re-rankings -- Monero/Aptos swap to rows 37/38, OKB/Mantle swap to
# rows 42/43) described by the commit "Updated cryptos list ... with GitHub
# Actions". Only the cells that actually changed are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''61.973.24'
$ws.Range("E2").Value = '  -1.51%  '

# Row 3
$ws.Range("D3").Value = '''3.412.02'
$ws.Range("E3").Value = '  -1.67%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '''576.14'
$ws.Range("E5").Value = '  -0.26%  '

# Row 6
$ws.Range("D6").Value = '''148.13'
$ws.Range("E6").Value = '  -0.01%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("E8").Value = '  +1.15%  '

# Row 9
$ws.Range("D9").Value = '''7.95'
$ws.Range("E9").Value = '  +4.05%  '

# Row 10
$ws.Range("E10").Value = '  -1.17%  '

# Row 11
$ws.Range("D11").Value = '''0.415'
$ws.Range("E11").Value = '  +2.95%  '

# Row 12
$ws.Range("D12").Value = '''3.993.35'
$ws.Range("E12").Value = '  -1.73%  '

# Row 13
$ws.Range("E13").Value = '  -0.09%  '

# Row 14
$ws.Range("D14").Value = '''28.40'
$ws.Range("E14").Value = '  -4.46%  '

# Row 15
$ws.Range("D15").Value = '''3.404.36'
$ws.Range("E15").Value = '  -1.84%  '

# Row 16
$ws.Range("E16").Value = '  +0.03%  '

# Row 17
$ws.Range("D17").Value = '''61.929.61'
$ws.Range("E17").Value = '  -1.63%  '

# Row 18
$ws.Range("D18").Value = '''6.40'
$ws.Range("E18").Value = '  +1.15%  '

# Row 19
$ws.Range("D19").Value = '''14.58'
$ws.Range("E19").Value = '  +1.36%  '

# Row 20
$ws.Range("D20").Value = '''8.96'
$ws.Range("E20").Value = '  -2.92%  '

# Row 21
$ws.Range("D21").Value = '''380.51'
$ws.Range("E21").Value = '  -1.96%  '

# Row 22
$ws.Range("E22").Value = '  +1.72%  '

# Row 23
$ws.Range("D23").Value = '''74.83'
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("E24").Value = '  +0.01%  '

# Row 25
$ws.Range("D25").Value = '''3.568.53'
$ws.Range("E25").Value = '  -1.12%  '

# Row 26
$ws.Range("E26").Value = '  -2.75%  '

# Row 27
$ws.Range("D27").Value = '''0.179'
$ws.Range("E27").Value = '  -0.22%  '

# Row 28
$ws.Range("D28").Value = '''7.60'
$ws.Range("E28").Value = '  +0.42%  '

# Row 29
$ws.Range("E29").Value = '  -0.01%  '

# Row 30
$ws.Range("D30").Value = '''7.91'
$ws.Range("E30").Value = '  -2.88%  '

# Row 31
$ws.Range("E31").Value = '  -0.02%  '

# Row 32
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("D33").Value = '''1.34'
$ws.Range("E33").Value = '  -1.68%  '

# Row 34
$ws.Range("D34").Value = '''23.05'
$ws.Range("E34").Value = '  -2.68%  '

# Row 35
$ws.Range("D35").Value = '''5.47'
$ws.Range("E35").Value = '  +3.80%  '

# Row 36
$ws.Range("E36").Value = '  +3.24%  '

# Row 37
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''169.77'
$ws.Range("E37").Value = '  -0.03%  '

# Row 38
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '''6.91'
$ws.Range("E38").Value = '  -2.10%  '

# Row 39
$ws.Range("D39").Value = '''30.20'
$ws.Range("E39").Value = '  -5.83%  '

# Row 40
$ws.Range("D40").Value = '''3.445.39'
$ws.Range("E40").Value = '  -1.76%  '

# Row 41
$ws.Range("D41").Value = '''0.0781'
$ws.Range("E41").Value = '  +3.53%  '

# Row 42
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''42.41'
$ws.Range("E42").Value = '  +0.16%  '

# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '''0.778'
$ws.Range("E43").Value = '  -2.76%  '

# Row 44
$ws.Range("E44").Value = '  -2.09%  '

# Row 45
$ws.Range("E45").Value = '  -2.32%  '

# Row 46
$ws.Range("E46").Value = '  -2.97%  '

# Row 47
$ws.Range("D47").Value = '''2.541.43'
$ws.Range("E47").Value = '  -2.80%  '

# Row 48
$ws.Range("D48").Value = '''6.90'

# Row 49
$ws.Range("D49").Value = '''22.74'
$ws.Range("E49").Value = '  -0.86%  '

# Row 50
$ws.Range("E50").Value = '  -4.35%  '

# Row 51
$ws.Range("D51").Value = '''0.999'
$ws.Range("E51").Value = '  -0.18%  '
